# Update segmentation with new timing:
# The "Participant Info" block (previously in columns I:L) is moved three
# columns to the right, now occupying columns L:O. This mirrors a manual
# re-segmentation of the worksheet where 3 blank columns were inserted
# before the block.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("In Lab")

# Source block containing the "Participant Info" table (rows 1-15).
$src = $ws.Range("I1:L15")

# Destination block, shifted right by 3 columns (I->L, J->M, K->N, L->O).
$dst = $ws.Range("L1:O15")

# Copy source (values + formatting) onto the destination, then clear the
# now-vacated source cells so nothing is duplicated.
$src.Copy($dst)
$src.Clear()

# Update column widths for the block's new location (was cols I,J,K -> now L,M,N).
$ws.Columns.Item(12).ColumnWidth = 37
$ws.Columns.Item(13).ColumnWidth = 16.7109375
$ws.Columns.Item(14).ColumnWidth = 14.140625

# Restore default width for the vacated columns I,J,K (no longer custom).
$ws.Columns.Item(9).ColumnWidth = 8.43
$ws.Columns.Item(10).ColumnWidth = 8.43
$ws.Columns.Item(11).ColumnWidth = 8.43

# Update the view: scrolled to show column B, with M15 selected.
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("M15").Select()
